$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

# Row 2
Set-TextValue 2 4 "36.923.83"
Set-TextValue 2 5 "  -0.24%  "

# Row 3
Set-TextValue 3 4 "2.039.71"
Set-TextValue 3 5 "  -0.33%  "

# Row 4
Set-TextValue 4 5 "  +0.20%  "

# Row 5
Set-TextValue 5 4 "244.48"
Set-TextValue 5 5 "  -1.32%  "

# Row 6
Set-TextValue 6 4 "0.655"
Set-TextValue 6 5 "  -0.90%  "

# Row 7
Set-TextValue 7 2 "USDC"
Set-TextValue 7 3 "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
Set-TextValue 7 4 "1.00"
Set-TextValue 7 5 "  -0.02%  "

# Row 8
Set-TextValue 8 2 "Solana"
Set-TextValue 8 3 "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
Set-TextValue 8 4 "57.11"
Set-TextValue 8 5 "  -0.82%  "

# Row 9
Set-TextValue 9 5 "  -0.90%  "

# Row 10
Set-TextValue 10 4 "0.0764"
Set-TextValue 10 5 "  -1.19%  "

# Row 11
Set-TextValue 11 5 "  +1.97%  "

# Row 12
Set-TextValue 12 4 "15.34"
Set-TextValue 12 5 "  -2.00%  "

# Row 13
Set-TextValue 13 4 "0.875"
Set-TextValue 13 5 "  +8.94%  "

# Row 14
Set-TextValue 14 4 "2.339.78"
Set-TextValue 14 5 "  -0.26%  "

# Row 15
Set-TextValue 15 4 "5.60"
Set-TextValue 15 5 "  +1.73%  "

# Row 16
Set-TextValue 16 4 "2.009.87"
Set-TextValue 16 5 "  -1.79%  "

# Row 17
Set-TextValue 17 4 "18.02"
Set-TextValue 17 5 "  +9.92%  "

# Row 18
Set-TextValue 18 4 "36.886.27"
Set-TextValue 18 5 "  -0.19%  "

# Row 19
Set-TextValue 19 4 "73.40"
Set-TextValue 19 5 "  -1.37%  "

# Row 20
Set-TextValue 20 4 "0.0₃0881"
Set-TextValue 20 5 "  -1.17%  "

# Row 21
Set-TextValue 21 4 "5.34"
Set-TextValue 21 5 "  +0.60%  "

# Row 22
Set-TextValue 22 4 "234.95"
Set-TextValue 22 5 "  -0.28%  "

# Row 23
Set-TextValue 23 4 "0.999"
Set-TextValue 23 5 "  -0.03%  "

# Row 24
Set-TextValue 24 4 "2.44"
Set-TextValue 24 5 "  +2.99%  "

# Row 25
Set-TextValue 25 4 "9.62"
Set-TextValue 25 5 "  +6.04%  "

# Row 26
Set-TextValue 26 4 "169.08"
Set-TextValue 26 5 "  +1.19%  "

# Row 27
Set-TextValue 27 4 "2.09"
Set-TextValue 27 5 "  -3.94%  "

# Row 28
Set-TextValue 28 4 "19.80"
Set-TextValue 28 5 "  +0.67%  "

# Row 29
Set-TextValue 29 5 "  +15.42%  "

# Row 30
Set-TextValue 30 5 "  -0.55%  "

# Row 31
Set-TextValue 31 5 "  -3.20%  "

# Row 32
Set-TextValue 32 4 "4.63"
Set-TextValue 32 5 "  +4.95%  "

# Row 33
Set-TextValue 33 4 "0.0608"
Set-TextValue 33 5 "  -0.21%  "

# Row 34
Set-TextValue 34 5 "  +0.02%  "

# Row 35
Set-TextValue 35 4 "0.0866"
Set-TextValue 35 5 "  -2.84%  "

# Row 36
Set-TextValue 36 5 "  +6.81%  "

# Row 37
Set-TextValue 37 4 "2.22"
Set-TextValue 37 5 "  +0.86%  "

# Row 38
Set-TextValue 38 5 "  -3.16%  "

# Row 39
Set-TextValue 39 5 "  -1.78%  "

# Row 40
Set-TextValue 40 5 "  +2.88%  "

# Row 41
Set-TextValue 41 4 "0.0974"
Set-TextValue 41 5 "  -8.00%  "

# Row 42
Set-TextValue 42 4 "0.0220"
Set-TextValue 42 5 "  +0.24%  "

# Row 43
Set-TextValue 43 5 "  +0.99%  "

# Row 44
Set-TextValue 44 4 "96.43"
Set-TextValue 44 5 "  +1.50%  "

# Row 45
Set-TextValue 45 4 "16.76"
Set-TextValue 45 5 "  -2.74%  "

# Row 46
Set-TextValue 46 4 "3.86"
Set-TextValue 46 5 "  +15.23%  "

# Row 47
Set-TextValue 47 4 "1.286.37"
Set-TextValue 47 5 "  +0.98%  "

# Row 48
Set-TextValue 48 4 "2.32"
Set-TextValue 48 5 "  -3.65%  "

# Row 49
Set-TextValue 49 5 "  -0.09%  "

# Row 50
Set-TextValue 50 2 "RocketPoolETH"
Set-TextValue 50 3 "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
Set-TextValue 50 4 "2.224.91"
Set-TextValue 50 5 "  -0.28%  "

# Row 51
Set-TextValue 51 2 "FraxShare"
Set-TextValue 51 3 "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue 51 4 "6.67"
Set-TextValue 51 5 "  +0.55%  "
